# Apply the updates described by the commit diff to the "Toiletten" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Toiletten")

# 1) Update the opening-hours text in column D (rows 5-19): switch the
#    "HH:MM–HH:MM" time notation to the shorter "H–H Uhr" notation.
$ws.Range("D5").Value  = "Badstr. 1-3 (0–24 Uhr)"
$ws.Range("D6").Value  = "Badstr. 4 (6–20 Uhr)"
$ws.Range("D7").Value  = "Hardenbergplatz 9-11 (6–22 Uhr)"
$ws.Range("D8").Value  = "Georgenstraße 14 (6–22 Uhr)"
$ws.Range("D9").Value  = "Dircksenstraße 2 (6–22 Uhr)"
$ws.Range("D10").Value = "Kronprinzessinnenweg 120 (11–21 Uhr)`n"
$ws.Range("D11").Value = "Sonntagstr. 37 (6–23 Uhr)"
$ws.Range("D12").Value = "Weitlingstraße 22 (6–22 Uhr)"
$ws.Range("D13").Value = "Kronprinzessinnenweg 250 (0–24 Uhr)"
$ws.Range("D14").Value = "Seegefelder Str. 1 (6–22 Uhr)"
$ws.Range("D15").Value = "General-Pape-Straße 1 (0–24 Uhr)"
$ws.Range("D16").Value = "Hildegard-Knef-Platz (0–24 Uhr)"
$ws.Range("D17").Value = "Koppenstraße 3 (6–22 Uhr)"
$ws.Range("D18").Value = "Potsdamer Platz 1 (0–24 Uhr)"
$ws.Range("D19").Value = "Europaplatz 1 (0–24 Uhr)"

# 2) The shorter text re-wraps differently, so a handful of rows now need
#    an explicit row height (column D has wrapText enabled).
$ws.Rows.Item(7).RowHeight  = 28.8
$ws.Rows.Item(8).RowHeight  = 28.8
$ws.Rows.Item(9).RowHeight  = 28.8
$ws.Rows.Item(10).RowHeight = 43.2
$ws.Rows.Item(12).RowHeight = 28.8
$ws.Rows.Item(13).RowHeight = 28.8
$ws.Rows.Item(18).RowHeight = 28.8

# 3) Move the active selection on the sheet to E18.
$ws.Range("E18").Select()
